$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("64000data")
Write-Host $ws.Name
